# relacionales.xlsx UI touch-up:
#  - rename the "usuarios registrados" sheet to "usuarios_registrados"
#  - move the selection/active-cell on "usuarios registrados" (now
#    usuarios_registrados) from B4 to F24, and it is no longer the active tab
#  - make "categorias_id" the active tab, with the selection/active cell
#    moved from C2 to E5

$wb = $excel.ActiveWorkbook

$wsCategorias = $wb.Worksheets.Item("categorias_id")
$wsUsuarios   = $wb.Worksheets.Item("usuarios registrados")

# Rename the third sheet (space -> underscore)
$wsUsuarios.Name = "usuarios_registrados"

# Update the (now inactive) usuarios_registrados sheet's selection first
$wsUsuarios.Activate()
$wsUsuarios.Range("F24").Select()

# Finish with categorias_id active/selected, becoming the workbook's active tab
$wsCategorias.Activate()
$wsCategorias.Range("E5").Select()
